# Weekly fruit/vegetable price update.
# Insert two new rows of data (one new reporting date: 44491) right before
# the existing row 132, pushing the rest of the table down by two rows
# (old rows 132-148 become 134-150).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 132, shifting rows 132:148 down to 134:150.
$ws.Range("A132:A133").EntireRow.Insert()

# --- New row 132 ---
$ws.Cells.Item(132, 1).Value = 4
$ws.Cells.Item(132, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(132, 3).Value = "Los Lagos"
$ws.Cells.Item(132, 4).Value = 44491
$ws.Cells.Item(132, 5).Value = 10
$ws.Cells.Item(132, 6).Value = 100112017
$ws.Cells.Item(132, 7).Value = "Apio"
$ws.Cells.Item(132, 8).Value = "Americana (o)"
$ws.Cells.Item(132, 9).Value = "Primera"
$ws.Cells.Item(132, 10).Value = 25
$ws.Cells.Item(132, 11).Value = 11000
$ws.Cells.Item(132, 12).Value = 11000
$ws.Cells.Item(132, 13).Value = 11000
$ws.Cells.Item(132, 14).Value = "`$/docena de matas"
$ws.Cells.Item(132, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(132, 16).Value = 1833
$ws.Cells.Item(132, 17).Value = 6
$ws.Cells.Item(132, 18).Value = "Hortaliza"

# --- New row 133 ---
$ws.Cells.Item(133, 1).Value = 4
$ws.Cells.Item(133, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(133, 3).Value = "Los Lagos"
$ws.Cells.Item(133, 4).Value = 44491
$ws.Cells.Item(133, 5).Value = 10
$ws.Cells.Item(133, 6).Value = 100112017
$ws.Cells.Item(133, 7).Value = "Apio"
$ws.Cells.Item(133, 8).Value = "Americana (o)"
$ws.Cells.Item(133, 9).Value = "Segunda"
$ws.Cells.Item(133, 10).Value = 25
$ws.Cells.Item(133, 11).Value = 10000
$ws.Cells.Item(133, 12).Value = 10000
$ws.Cells.Item(133, 13).Value = 10000
$ws.Cells.Item(133, 14).Value = "`$/docena de matas"
$ws.Cells.Item(133, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(133, 16).Value = 1667
$ws.Cells.Item(133, 17).Value = 6
$ws.Cells.Item(133, 18).Value = "Hortaliza"

# Make sure the date cells use the same date/time number format as the
# rest of the "Fecha" column (style index carried over from the insert,
# but set explicitly to be safe).
$ws.Range("D132:D133").NumberFormat = "YYYY-MM-DD HH:MM:SS"
